# Updates the cryptos list (Price / Volume(1h) columns, and a couple of
# coin name/link/price swaps) to match the latest scrape.
#
# All cells in this sheet are stored as plain text (inlineStr) in the
# source workbook - including values that look numeric, like "244.89" or
# "1.00" - so that things like trailing zeros and thousands-separated
# "prices" (e.g. "97.475.64") are preserved exactly as scraped. Excel's
# COM Value setter auto-detects numeric-looking strings and would coerce
# them into real floating point numbers (losing formatting and precision,
# e.g. "33.03" -> 33.030000000000001). To avoid that, every write below
# temporarily forces the cell's number format to Text ("@") before
# assigning the value, then restores the cell style to "Normal" so the
# on-disk style (no explicit style index) matches the original.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($address, $value) {
    $rng = $ws.Range($address)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# --- Row 2 : Bitcoin ---
Set-TextValue "D2" "97.480.60"
Set-TextValue "E2" "  +2.02%  "

# --- Row 3 : Ethereum ---
Set-TextValue "D3" "3.599.32"
Set-TextValue "E3" "  +0.21%  "

# --- Row 4 : TetherUSD ---
Set-TextValue "E4" "  +0.01%  "

# --- Row 5 : Solana ---
Set-TextValue "D5" "244.58"
Set-TextValue "E5" "  +3.04%  "

# --- Row 6 : XRP ---
Set-TextValue "E6" "  +17.39%  "

# --- Row 7 ---
Set-TextValue "D7" "653.27"
Set-TextValue "E7" "  -0.86%  "

# --- Row 8 ---
Set-TextValue "D8" "0.428"
Set-TextValue "E8" "  +6.19%  "

# --- Row 9 ---
Set-TextValue "E9" "  -0.08%  "

# --- Row 10 ---
Set-TextValue "E10" "  +2.28%  "

# --- Row 11 ---
Set-TextValue "D11" "3.596.91"
Set-TextValue "E11" "  +0.24%  "

# --- Row 12 ---
Set-TextValue "D12" "44.76"
Set-TextValue "E12" "  +4.32%  "

# --- Row 13 ---
Set-TextValue "E13" "  +0.91%  "

# --- Row 14 ---
Set-TextValue "D14" "6.49"
Set-TextValue "E14" "  -0.07%  "

# --- Row 15 ---
Set-TextValue "D15" "4.265.94"
Set-TextValue "E15" "  +0.14%  "

# --- Row 16 ---
Set-TextValue "D16" "97.345.72"
Set-TextValue "E16" "  +2.03%  "

# --- Row 17 ---
Set-TextValue "E17" "  +2.66%  "

# --- Row 18 ---
Set-TextValue "D18" "3.604.84"
Set-TextValue "E18" "  +0.35%  "

# --- Row 19 ---
Set-TextValue "E19" "  +0.27%  "

# --- Row 20 ---
Set-TextValue "E20" "  -1.20%  "

# --- Row 21 ---
Set-TextValue "D21" "18.24"
Set-TextValue "E21" "  +1.28%  "

# --- Row 22 ---
Set-TextValue "D22" "0.528"
Set-TextValue "E22" "  +8.09%  "

# --- Row 23 / 24 : SuiNetwork and BitcoinCash swap places ---
Set-TextValue "B23" "BitcoinCash"
Set-TextValue "C23" "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D23" "520.00"
Set-TextValue "E23" "  +1.74%  "

Set-TextValue "B24" "SuiNetwork"
Set-TextValue "C24" "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue "D24" "3.49"
Set-TextValue "E24" "  +0.27%  "

# --- Row 25 ---
Set-TextValue "E25" "  +4.67%  "

# --- Row 26 ---
Set-TextValue "D26" "6.99"
Set-TextValue "E26" "  -2.16%  "

# --- Row 27 ---
Set-TextValue "D27" "103.96"
Set-TextValue "E27" "  +8.41%  "

# --- Row 28 ---
Set-TextValue "E28" "  +3.01%  "

# --- Row 29 ---
Set-TextValue "E29" "  +22.56%  "

# --- Row 30 ---
Set-TextValue "D30" "3.793.65"
Set-TextValue "E30" "  +0.26%  "

# --- Row 31 ---
Set-TextValue "E31" "  -2.59%  "

# --- Row 32 ---
Set-TextValue "E32" "  +3.38%  "

# --- Row 33 ---
Set-TextValue "D33" "0.999"
Set-TextValue "E33" "  -0.40%  "

# --- Row 34 ---
Set-TextValue "E34" "  +5.98%  "

# --- Row 35 ---
Set-TextValue "D35" "0.989"
Set-TextValue "E35" "  -0.76%  "

# --- Row 36 ---
Set-TextValue "D36" "31.92"

# --- Row 37 ---
Set-TextValue "E37" "  +3.62%  "

# --- Row 38 ---
Set-TextValue "E38" "  -1.43%  "

# --- Row 39 ---
Set-TextValue "D39" "618.79"
Set-TextValue "E39" "  +2.64%  "

# --- Row 40 ---
Set-TextValue "D40" "8.80"
Set-TextValue "E40" "  +0.68%  "

# --- Row 41 ---
Set-TextValue "E41" "  +2.03%  "

# --- Row 42 ---
Set-TextValue "D42" "1.94"
Set-TextValue "E42" "  +2.31%  "

# --- Row 43 ---
Set-TextValue "E43" "  +1.84%  "

# --- Row 44 ---
Set-TextValue "E44" "  -0.09%  "

# --- Row 45 ---
Set-TextValue "E45" "  +6.46%  "

# --- Row 46 ---
Set-TextValue "D46" "0.448"
Set-TextValue "E46" "  +42.86%  "

# --- Row 47 ---
Set-TextValue "E47" "  +6.36%  "

# --- Row 48 ---
Set-TextValue "E48" "  +1.19%  "

# --- Row 49 ---
Set-TextValue "E49" "  +1.00%  "

# --- Row 50 ---
Set-TextValue "E50" "  +5.16%  "

# --- Row 51 : dogwifhat -> EnergySwap ---
Set-TextValue "B51" "EnergySwap"
Set-TextValue "C51" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D51" "33.03"
Set-TextValue "E51" "  -5.71%  "
